$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A192").Value = 191
$ws.Range("B192").Value = 1
$ws.Range("C192").Value = "2024-06-19 00:57:52"
$ws.Range("D192").Value = 200
$ws.Range("E192").Value = 15

$ws.Range("A193").Value = 192
$ws.Range("B193").Value = 2
$ws.Range("C193").Value = "2024-06-19 00:57:53"
$ws.Range("D193").Value = 200
$ws.Range("E193").Value = 3
